# "Raw and Clean data from SSA for June 25th"
#
# - Row 25 col A currently holds the text "2020-06-24" (a shared string).
#   It becomes a real date serial (2020-06-24 == 44006) formatted with the
#   workbook's existing "YYYY-MM-DD HH:MM:SS" date/time number format.
# - A brand-new row 26 is appended with the June 25th data, whose date cell
#   uses a new "YYYY-MM-DD" date-only number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: replace the text date with a true date value, reusing the
# existing datetime number format already used by A2:A24.
$ws.Range("A25").Value = 44006
$ws.Range("A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 26: newly added data for 2020-06-25.
$ws.Range("A26").Value = 44007
$ws.Range("A26").NumberFormat = "YYYY-MM-DD"

$ws.Range("B26").Value = 202951
$ws.Range("C26").Value = 262117
$ws.Range("D26").Value = 63583
$ws.Range("E26").Value = 25060
$ws.Range("F26").Value = 31.46
